# Generate Report for Archive
# - Swap the two rows for "6958c7be...md" and "85991453...png" so the
#   ".md" entry sorts before the ".png" entry (rows 4 & 5) on every sheet.
# - Flip every remaining "Ready for handoff" status to "In Translation"
#   (the report now reflects the archived/processed state).

$wb = $excel.ActiveWorkbook

# ---- Overview sheet ----------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A4").Value = "6958c7be-b92b-49cc-bc67-a852aff2c471.md"
$ws.Range("B4").Value = "e2e\6958c7be-b92b-49cc-bc67-a852aff2c471.md"
$ws.Range("C4").Value = ".md"
$ws.Range("E4").Value = "In Translation"
$ws.Range("F4").Value = "In Translation"
$ws.Range("G4").Value = "2016-08-22 04:43:50"

$ws.Range("A5").Value = "85991453-56d7-4a31-b418-976de62e35ba.png"
$ws.Range("B5").Value = "e2e\85991453-56d7-4a31-b418-976de62e35ba.png"
$ws.Range("C5").Value = ".png"
$ws.Range("E5").Value = "In Translation"
$ws.Range("F5").Value = "In Translation"
$ws.Range("G5").Value = "2016-08-22 04:43:12"

$ws.Range("E6").Value = "In Translation"
$ws.Range("F6").Value = "In Translation"

$ws.Range("E7").Value = "In Translation"
$ws.Range("F7").Value = "In Translation"

# ---- zh-cn sheet ---------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A4").Value = "6958c7be-b92b-49cc-bc67-a852aff2c471.md"
$ws.Range("B4").Value = ".md"
$ws.Range("C4").Value = "In Translation"
$ws.Range("G4").Value = "6958c7be-b92b-49cc-bc67-a852aff2c471.67e8468d284d4311ea8d3e11e125777329db808c.zh-cn.xlf"
$ws.Range("H4").Value = "2016-08-22 04:43:45"
$ws.Range("M4").Value = "'True"
$ws.Range("N4").Value = "'"

$ws.Range("A5").Value = "85991453-56d7-4a31-b418-976de62e35ba.png"
$ws.Range("B5").Value = ".png"
$ws.Range("C5").Value = "In Translation"
$ws.Range("G5").Value = "fa61a0e5f771c450c9ad4eb8d72f33063fae7f41.png"
$ws.Range("H5").Value = "2016-08-22 04:43:08"
$ws.Range("M5").Value = "True(Dependency)"
$ws.Range("N5").Value = "e2e\0b98d60b-d99f-4c3e-9575-704f305e38a7.md"

$ws.Range("C6").Value = "In Translation"
$ws.Range("C7").Value = "In Translation"

# ---- de-de sheet ---------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A4").Value = "6958c7be-b92b-49cc-bc67-a852aff2c471.md"
$ws.Range("B4").Value = ".md"
$ws.Range("C4").Value = "In Translation"
$ws.Range("G4").Value = "6958c7be-b92b-49cc-bc67-a852aff2c471.67e8468d284d4311ea8d3e11e125777329db808c.de-de.xlf"
$ws.Range("H4").Value = "2016-08-22 04:43:50"
$ws.Range("M4").Value = "'True"
$ws.Range("N4").Value = "'"

$ws.Range("A5").Value = "85991453-56d7-4a31-b418-976de62e35ba.png"
$ws.Range("B5").Value = ".png"
$ws.Range("C5").Value = "In Translation"
$ws.Range("G5").Value = "fa61a0e5f771c450c9ad4eb8d72f33063fae7f41.png"
$ws.Range("H5").Value = "2016-08-22 04:43:12"
$ws.Range("M5").Value = "True(Dependency)"
$ws.Range("N5").Value = "e2e\0b98d60b-d99f-4c3e-9575-704f305e38a7.md"

$ws.Range("C6").Value = "In Translation"
$ws.Range("C7").Value = "In Translation"
